$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C51")
$range.Value = 0.01
$range.NumberFormat = "#,##0.0000"
